$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Python Fundamentals"
$ws.Range("C3").Value = "Brown"
$ws.Range("D3").Value = "T112"
$ws.Range("E3").Value = 1
